# Weekly update: insert 4 new "Chirimoya" price rows (week of 2023-11-28)
# above the existing data block, pushing the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at the top of the data block (old row 327 -> new row 331, etc.)
$ws.Rows("327:330").Insert()

# Columns that are constant across this whole product block
$ws.Range("A327:A330").Value = 8
$ws.Range("B327:B330").Value = "Terminal La Palmera de La Serena"
$ws.Range("C327:C330").Value = "Coquimbo"
$ws.Range("D327:D330").Value = 45258
$ws.Range("E327:E330").Value = 4
$ws.Range("F327:F330").Value = "Fruta"
$ws.Range("G327:G330").Value = 100107
$ws.Range("H327:H330").Value = "Otros"
$ws.Range("I327:I330").Value = 100107002
$ws.Range("J327:J330").Value = "Chirimoya"
$ws.Range("K327:K330").Value = "Cultivar IV Región"
$ws.Range("Q327:Q330").Value = "`$/bandeja 10 kilos"
$ws.Range("R327:R330").Value = "Provincia de Limarí"
$ws.Range("T327:T330").Value = 10

# Row 327 - Especial
$ws.Range("L327").Value = "Especial"
$ws.Range("M327").Value = 480
$ws.Range("N327").Value = 17000
$ws.Range("O327").Value = 18000
$ws.Range("P327").Value = 17500
$ws.Range("S327").Value = 1750

# Row 328 - Primera
$ws.Range("L328").Value = "Primera"
$ws.Range("M328").Value = 360
$ws.Range("N328").Value = 14000
$ws.Range("O328").Value = 15000
$ws.Range("P328").Value = 14500
$ws.Range("S328").Value = 1450

# Row 329 - Segunda
$ws.Range("L329").Value = "Segunda"
$ws.Range("M329").Value = 240
$ws.Range("N329").Value = 10000
$ws.Range("O329").Value = 11000
$ws.Range("P329").Value = 10500
$ws.Range("S329").Value = 1050

# Row 330 - Tercera
$ws.Range("L330").Value = "Tercera"
$ws.Range("M330").Value = 200
$ws.Range("N330").Value = 7000
$ws.Range("O330").Value = 8000
$ws.Range("P330").Value = 7500
$ws.Range("S330").Value = 750
